$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D5").Value = 44524
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103003
$ws.Range("J5").Value = "Damasco"
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 27000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 27500
$ws.Range("Q5").Value = "$/bandeja 18 kilos"
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 1528
$ws.Range("T5").Value = 18
